{"js": "// Apply the benchmark-stats corrections to the single results table.\n// Row numbers below are 1-based (as seen in the document); table rows\n// are addressed with the 0-based `getCell(rowIndex, 0)` API.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of 0-based row index -> new cell text.\nconst updates = {\n  0: \"0M\",       // row 1: 99.98 -> 0M\n  1: \"0M\",       // row 2: 0.65 -> 0M\n  2: \"0M\",       // row 3: 3844 -> 0M\n  3: \"3636\",     // row 4: 909 -> 3636\n  5: \"0.00070\",  // row 6: 0.00030 -> 0.00070\n  6: \"0.00018\",  // row 7: 0.00009 -> 0.00018\n  7: \"0.00005\",  // row 8: 0.00002 -> 0.00005\n  8: \"0.00027\",  // row 9: 0.00008 -> 0.00027\n  9: \"0.00033\",  // row 10: 0.00009 -> 0.00033\n  10: \"0.00037\", // row 11: 0.00010 -> 0.00037\n  11: \"0.65493\", // row 12: 0.08352 -> 0.65493\n  43: \"99.98\",   // row 44: collapsed multi-run row -> 99.98\n  44: \"0.65\",    // row 45: collapsed multi-run row -> 0.65\n  45: \"3844\",    // row 46: collapsed multi-run row -> 3844\n};\n\nfor (const rowIndex of Object.keys(updates)) {\n  const cell = table.getCell(Number(rowIndex), 0);\n  cell.value = updates[rowIndex];\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-stats corrections to the single results table.\n# Rows are addressed with the COM 1-based Rows.Item(n)/Cells.Item(1) API.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @{\n  1  = \"0M\"       # row 1: 99.98 -> 0M\n  2  = \"0M\"       # row 2: 0.65 -> 0M\n  3  = \"0M\"       # row 3: 3844 -> 0M\n  4  = \"3636\"     # row 4: 909 -> 3636\n  6  = \"0.00070\"  # row 6: 0.00030 -> 0.00070\n  7  = \"0.00018\"  # row 7: 0.00009 -> 0.00018\n  8  = \"0.00005\"  # row 8: 0.00002 -> 0.00005\n  9  = \"0.00027\"  # row 9: 0.00008 -> 0.00027\n  10 = \"0.00033\"  # row 10: 0.00009 -> 0.00033\n  11 = \"0.00037\"  # row 11: 0.00010 -> 0.00037\n  12 = \"0.65493\"  # row 12: 0.08352 -> 0.65493\n  44 = \"99.98\"    # row 44: collapsed multi-run row -> 99.98\n  45 = \"0.65\"     # row 45: collapsed multi-run row -> 0.65\n  46 = \"3844\"     # row 46: collapsed multi-run row -> 3844\n}\n\nforeach ($rowNum in $updates.Keys) {\n  $cell = $t.Rows.Item($rowNum).Cells.Item(1)\n  $cell.Range.Text = $updates[$rowNum]\n}\n"}
